$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Bitcoin"
$ws.Range("C2").Value = "https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc"
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.002.45"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -2.54%  "

$ws.Range("B3").Value = "Ethereum"
$ws.Range("C3").Value = "https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.909.61"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -4.84%  "

$ws.Range("B4").Value = "TetherUSD"
$ws.Range("C4").Value = "https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("B5").Value = "BNB"
$ws.Range("C5").Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "323.64"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -1.68%  "

$ws.Range("B6").Value = "USDC"
$ws.Range("C6").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.002"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.14%  "

$ws.Range("B7").Value = "XRP"
$ws.Range("C7").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4599"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -1.89%  "

$ws.Range("B8").Value = "Cardano"
$ws.Range("C8").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3819"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -3.59%  "

$ws.Range("B9").Value = "Dogecoin"
$ws.Range("C9").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07715"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -4.19%  "

$ws.Range("B10").Value = "Polygon"
$ws.Range("C10").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9800"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -2.17%  "

$ws.Range("B11").Value = "Solana"
$ws.Range("C11").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "21.99"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -4.71%  "

$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.915.69"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -5.36%  "

$ws.Range("B13").Value = "Chainlink"
$ws.Range("C13").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.955"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -4.21%  "

$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.667"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -3.73%  "

$ws.Range("B15").Value = "TRON"
$ws.Range("C15").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.07054"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -0.71%  "

$ws.Range("B16").Value = "BinanceUSD"
$ws.Range("C16").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.004"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -0.07%  "

$ws.Range("B17").Value = "Litecoin"
$ws.Range("C17").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "83.75"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -6.05%  "

$ws.Range("B18").Value = "ShibaInu"
$ws.Range("C18").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000009519"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -5.43%  "

$ws.Range("B19").Value = "Avalanche"
$ws.Range("C19").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "16.69"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -4.39%  "

$ws.Range("B20").Value = "Dai"
$ws.Range("C20").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.002"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +0.00%  "

$ws.Range("B21").Value = "WrappedBTC"
$ws.Range("C21").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "29.007.45"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -2.62%  "

$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.317"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -5.05%  "

$ws.Range("B23").Value = "Cosmos"
$ws.Range("C23").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.89"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -3.67%  "

$ws.Range("B24").Value = "Toncoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.097"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -0.51%  "

$ws.Range("B25").Value = "Monero"
$ws.Range("C25").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "157.59"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -1.19%  "

$ws.Range("B26").Value = "EthereumClassic"
$ws.Range("C26").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "19.10"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -3.14%  "

$ws.Range("B27").Value = "InternetComputer(DFINITY)"
$ws.Range("C27").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "5.617"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -5.58%  "

$ws.Range("B28").Value = "BitcoinCash"
$ws.Range("C28").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "117.55"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -2.47%  "

$ws.Range("B29").Value = "LidoDAOToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.853"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -4.38%  "

$ws.Range("B30").Value = "Stellar"
$ws.Range("C30").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.09266"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -2.35%  "

$ws.Range("B31").Value = "ImmutableX"
$ws.Range("C31").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.8607"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -6.10%  "

$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.092"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -3.79%  "

$ws.Range("B33").Value = "ARBITRUM"
$ws.Range("C33").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.247"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -8.27%  "

$ws.Range("B34").Value = "HuobiToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.011"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -6.86%  "

$ws.Range("B35").Value = "Hedera"
$ws.Range("C35").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.05696"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -2.76%  "

$ws.Range("B36").Value = "TrustWalletToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.140"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -3.25%  "

$ws.Range("B37").Value = "Frax"
$ws.Range("C37").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.002"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -0.11%  "

$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02036"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -4.33%  "

$ws.Range("B39").Value = "FraxShare"
$ws.Range("C39").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "7.421"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -7.13%  "

$ws.Range("B40").Value = "TheSandbox"
$ws.Range("C40").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5496"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -5.14%  "

$ws.Range("B41").Value = "Algorand"
$ws.Range("C41").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1757"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -3.90%  "

$ws.Range("B42").Value = "Aptos"
$ws.Range("C42").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "9.297"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -5.95%  "

$ws.Range("B43").Value = "MXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.763"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -0.60%  "

$ws.Range("B44").Value = "Decentraland"
$ws.Range("C44").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5196"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -4.14%  "

$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "11.30"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -6.40%  "

$ws.Range("B46").Value = "RenderToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.076"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -5.22%  "

$ws.Range("B47").Value = "Cronos"
$ws.Range("C47").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.06794"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -2.79%  "

$ws.Range("B48").Value = "Quant"
$ws.Range("C48").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "111.45"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -3.14%  "

$ws.Range("B49").Value = "NEARProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.779"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -5.17%  "

$ws.Range("B50").Value = "PEPE"
$ws.Range("C50").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.000002502"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -29.36%  "

$ws.Range("B51").Value = "PaxDollar"
$ws.Range("C51").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.001"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -0.13%  "
